$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet "Phaser Links" (sheet2): add the new link rows ---
$ws2.Range("A2").Value = "https://www.codeandweb.com/texturepacker/tutorials/creating-spritesheets-for-phaser-with-texturepacker"
$ws2.Range("B2").Value = "Texture Packer"

$ws2.Range("A3").Value = "http://mightyfingers.com/"
$ws2.Range("B3").Value = "Phaser framework/editor"

$ws2.Range("A4").Value = "https://vimeo.com/album/3156158"
$ws2.Range("B4").Value = "mightyfingers videos"

$ws2.Range("A5").Value = "http://inkubator.io/inkubatethis.html"
$ws2.Range("B5").Value = "mightyfingers blog"

$ws2.Range("A6").Value = "https://elmvids.groob.io/"
$ws2.Range("B6").Value = "Elm Videos"

$ws2.Range("A7").Value = "https://www.youtube.com/watch?v=ZwaomOYGuYo"
$ws2.Range("B7").Value = "Tile video"

# Turn the Texture Packer URL into a real hyperlink (adds the Hyperlink style/font too)
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://www.codeandweb.com/texturepacker/tutorials/creating-spritesheets-for-phaser-with-texturepacker") | Out-Null

# Widen column A to fit the new long text (closest this engine's column-width
# rounding can land to the authored 61.36328125 OOXML width)
$ws2.Columns.Item(1).ColumnWidth = 60.5

# --- View state: Phaser Links becomes the active/selected sheet ---
$ws2.Range("B10").Select() | Out-Null
$ws2.Activate()
